# Apply updated power-flow (vm_pu) results for the 380 kV case.
# Slack bus voltage set-point changed from 1.05 to 1.02 p.u. (column B),
# which in turn changes every other bus voltage magnitude in the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$row2BF = New-Object "object[,]" 1,5
$row2BF[0,0] = 1.02
$row2BF[0,1] = 1.022816374890177
$row2BF[0,2] = 1.027559953783402
$row2BF[0,3] = 1.032988247842739
$row2BF[0,4] = 1.042879440361947
$ws.Range("B2:F2").Value = $row2BF
$row2IN = New-Object "object[,]" 1,6
$row2IN[0,0] = 1.02899998184734
$row2IN[0,1] = 1.028000085590942
$row2IN[0,2] = 1.030379327425508
$row2IN[0,3] = 1.03579189266329
$row2IN[0,4] = 1.045654881882074
$row2IN[0,5] = 1.013398075608869
$ws.Range("I2:N2").Value = $row2IN

# Row 3
$row3BF = New-Object "object[,]" 1,5
$row3BF[0,0] = 1.02
$row3BF[0,1] = 1.023725666834165
$row3BF[0,2] = 1.028231084591541
$row3BF[0,3] = 1.033869257034048
$row3BF[0,4] = 1.044013274956967
$ws.Range("B3:F3").Value = $row3BF
$row3IN = New-Object "object[,]" 1,6
$row3IN[0,0] = 1.029140952402609
$row3IN[0,1] = 1.028547797808808
$row3IN[0,2] = 1.030858578367409
$row3IN[0,3] = 1.036481595846159
$row3IN[0,4] = 1.04659879274947
$row3IN[0,5] = 1.013579247713899
$ws.Range("I3:N3").Value = $row3IN

# Row 4
$row4BF = New-Object "object[,]" 1,5
$row4BF[0,0] = 1.02
$row4BF[0,1] = 1.02431404750331
$row4BF[0,2] = 1.028664738682125
$row4BF[0,3] = 1.034439802239133
$row4BF[0,4] = 1.044747914831296
$ws.Range("B4:F4").Value = $row4BF
$row4IN = New-Object "object[,]" 1,6
$row4IN[0,0] = 1.029230095846343
$row4IN[0,1] = 1.028901606991762
$row4IN[0,2] = 1.031167396144163
$row4IN[0,3] = 1.036927718930909
$row4IN[0,4] = 1.047209956401147
$row4IN[0,5] = 1.013696260104713
$ws.Range("I4:N4").Value = $row4IN

# Row 5
$row5BF = New-Object "object[,]" 1,5
$row5BF[0,0] = 1.02
$row5BF[0,1] = 1.024561403770503
$row5BF[0,2] = 1.02884689877732
$row5BF[0,3] = 1.034679771723932
$row5BF[0,4] = 1.045056989683696
$ws.Range("B5:F5").Value = $row5BF
$row5IN = New-Object "object[,]" 1,6
$row5IN[0,0] = 1.029267074457312
$row5IN[0,1] = 1.029050204510481
$row5IN[0,2] = 1.031296913263204
$row5IN[0,3] = 1.037115229689299
$row5IN[0,4] = 1.047466982261458
$row5IN[0,5] = 1.013745399586935
$ws.Range("I5:N5").Value = $row5IN

# Row 6
$row6BF = New-Object "object[,]" 1,5
$row6BF[0,0] = 1.02
$row6BF[0,1] = 1.024602936011279
$row6BF[0,2] = 1.02887747553214
$row6BF[0,3] = 1.034720070210241
$row6BF[0,4] = 1.045108898310083
$ws.Range("B6:F6").Value = $row6BF
$row6IN = New-Object "object[,]" 1,6
$row6IN[0,0] = 1.029273254144216
$row6IN[0,1] = 1.029075146225924
$row6IN[0,2] = 1.031318641553042
$row6IN[0,3] = 1.037146711236469
$row6IN[0,4] = 1.047510143446122
$row6IN[0,5] = 1.013753647232708
$ws.Range("I6:N6").Value = $row6IN

# Row 7
$row7BF = New-Object "object[,]" 1,5
$row7BF[0,0] = 1.02
$row7BF[0,1] = 1.024317352686902
$row7BF[0,2] = 1.028667173297976
$row7BF[0,3] = 1.034443008283262
$row7BF[0,4] = 1.044752043792909
$ws.Range("B7:F7").Value = $row7BF
$row7IN = New-Object "object[,]" 1,6
$row7IN[0,0] = 1.029230591911381
$row7IN[0,1] = 1.028903593125476
$row7IN[0,2] = 1.031169127977289
$row7IN[0,3] = 1.036930224614671
$row7IN[0,4] = 1.047213390430129
$row7IN[0,5] = 1.01369691691568
$ws.Range("I7:N7").Value = $row7IN

# Row 8
$row8BF = New-Object "object[,]" 1,5
$row8BF[0,0] = 1.02
$row8BF[0,1] = 1.023123672182278
$row8BF[0,2] = 1.027786891507369
$row8BF[0,3] = 1.033285890390322
$row8BF[0,4] = 1.043262423165682
$ws.Range("B8:F8").Value = $row8BF
$row8IN = New-Object "object[,]" 1,6
$row8IN[0,0] = 1.029048052455761
$row8IN[0,1] = 1.028185310635995
$row8IN[0,2] = 1.030541558872346
$row8IN[0,3] = 1.036025013911484
$row8IN[0,4] = 1.045973799731858
$row8IN[0,5] = 1.013459348532804
$ws.Range("I8:N8").Value = $row8IN

# Row 9
$row9BF = New-Object "object[,]" 1,5
$row9BF[0,0] = 1.02
$row9BF[0,1] = 1.021020353938149
$row9BF[0,2] = 1.026231096395358
$row9BF[0,3] = 1.031250571599318
$row9BF[0,4] = 1.04064500466883
$ws.Range("B9:F9").Value = $row9BF
$row9IN = New-Object "object[,]" 1,6
$row9IN[0,0] = 1.028710545158183
$row9IN[0,1] = 1.026915069336979
$row9IN[0,2] = 1.029425874407922
$row9IN[0,3] = 1.034428720862707
$row9IN[0,4] = 1.043792503927648
$row9IN[0,5] = 1.01303906643594
$ws.Range("I9:N9").Value = $row9IN

# Row 10
$row10BF = New-Object "object[,]" 1,5
$row10BF[0,0] = 1.02
$row10BF[0,1] = 1.019618257935725
$row10BF[0,2] = 1.025190874445915
$row10BF[0,3] = 1.029896220100315
$row10BF[0,4] = 1.038905129510816
$ws.Range("B10:F10").Value = $row10BF
$row10IN = New-Object "object[,]" 1,6
$row10IN[0,0] = 1.028474929386494
$row10IN[0,1] = 1.026065242603252
$row10IN[0,2] = 1.028675542565752
$row10IN[0,3] = 1.033363769456626
$row10IN[0,4] = 1.042340380410411
$row10IN[0,5] = 1.01275778316008
$ws.Range("I10:N10").Value = $row10IN

# Row 11
$row11BF = New-Object "object[,]" 1,5
$row11BF[0,0] = 1.02
$row11BF[0,1] = 1.019011172172623
$row11BF[0,2] = 1.024739746804257
$row11BF[0,3] = 1.029310382538446
$row11BF[0,4] = 1.038152954494236
$ws.Range("B11:F11").Value = $row11BF
$row11IN = New-Object "object[,]" 1,6
$row11IN[0,0] = 1.028370397636827
$row11IN[0,1] = 1.02569655706187
$row11IN[0,2] = 1.028349102172932
$row11IN[0,3] = 1.032902464079779
$row11IN[0,4] = 1.041712094996036
$row11IN[0,5] = 1.012635728357418
$ws.Range("I11:N11").Value = $row11IN

# Row 12
$row12BF = New-Object "object[,]" 1,5
$row12BF[0,0] = 1.02
$row12BF[0,1] = 1.018785678860518
$row12BF[0,2] = 1.024572073398303
$row12BF[0,3] = 1.029092868279036
$row12BF[0,4] = 1.037873744370017
$ws.Range("B12:F12").Value = $row12BF
$row12IN = New-Object "object[,]" 1,6
$row12IN[0,0] = 1.028331193698861
$row12IN[0,1] = 1.025559505690137
$row12IN[0,2] = 1.028227617056236
$row12IN[0,3] = 1.032731089242572
$row12IN[0,4] = 1.041478796573182
$row12IN[0,5] = 1.012590353397531
$ws.Range("I12:N12").Value = $row12IN

# Row 13
$row13BF = New-Object "object[,]" 1,5
$row13BF[0,0] = 1.02
$row13BF[0,1] = 1.018834047687099
$row13BF[0,2] = 1.024608044620858
$row13BF[0,3] = 1.029139521652938
$row13BF[0,4] = 1.037933627682915
$ws.Range("B13:F13").Value = $row13BF
$row13IN = New-Object "object[,]" 1,6
$row13IN[0,0] = 1.02833962008908
$row13IN[0,1] = 1.025588908420438
$row13IN[0,2] = 1.028253686446064
$row13IN[0,3] = 1.032767850868434
$row13IN[0,4] = 1.041528836484595
$row13IN[0,5] = 1.012600088212333
$ws.Range("I13:N13").Value = $row13IN

# Row 14
$row14BF = New-Object "object[,]" 1,5
$row14BF[0,0] = 1.02
$row14BF[0,1] = 1.018992532711474
$row14BF[0,2] = 1.024725888992456
$row14BF[0,3] = 1.029292400861993
$row14BF[0,4] = 1.038129871199502
$ws.Range("B14:F14").Value = $row14BF
$row14IN = New-Object "object[,]" 1,6
$row14IN[0,0] = 1.028367164697311
$row14IN[0,1] = 1.025685230493264
$row14IN[0,2] = 1.028339064874054
$row14IN[0,3] = 1.032888298697924
$row14IN[0,4] = 1.041692808947256
$row14IN[0,5] = 1.012631978430221
$ws.Range("I14:N14").Value = $row14IN

# Row 15
$row15BF = New-Object "object[,]" 1,5
$row15BF[0,0] = 1.02
$row15BF[0,1] = 1.019090181298367
$row15BF[0,2] = 1.024798482990476
$row15BF[0,3] = 1.029386606992681
$row15BF[0,4] = 1.038250807328911
$ws.Range("B15:F15").Value = $row15BF
$row15IN = New-Object "object[,]" 1,6
$row15IN[0,0] = 1.028384086011022
$row15IN[0,1] = 1.025744563791324
$row15IN[0,2] = 1.02839163881207
$row15IN[0,3] = 1.032962507225984
$row15IN[0,4] = 1.041793847716559
$row15IN[0,5] = 1.012651621969781
$ws.Range("I15:N15").Value = $row15IN

# Row 16
$row16BF = New-Object "object[,]" 1,5
$row16BF[0,0] = 1.02
$row16BF[0,1] = 1.019658548933014
$row16BF[0,2] = 1.025220799571068
$row16BF[0,3] = 1.029935113046817
$row16BF[0,4] = 1.038955074267129
$ws.Range("B16:F16").Value = $row16BF
$row16IN = New-Object "object[,]" 1,6
$row16IN[0,0] = 1.028481814021074
$row16IN[0,1] = 1.026089696261254
$row16IN[0,2] = 1.028697174943051
$row16IN[0,3] = 1.033394381162511
$row16IN[0,4] = 1.04238208809512
$row16IN[0,5] = 1.012765878138786
$ws.Range("I16:N16").Value = $row16IN

# Row 17
$row17BF = New-Object "object[,]" 1,5
$row17BF[0,0] = 1.02
$row17BF[0,1] = 1.020015079747943
$row17BF[0,2] = 1.025485519904336
$row17BF[0,3] = 1.030279339179861
$row17BF[0,4] = 1.039397164596682
$ws.Range("B17:F17").Value = $row17BF
$row17IN = New-Object "object[,]" 1,6
$row17IN[0,0] = 1.028542445096334
$row17IN[0,1] = 1.026306000546719
$row17IN[0,2] = 1.028888417631373
$row17IN[0,3] = 1.033665238092887
$row17IN[0,4] = 1.042751208324283
$row17IN[0,5] = 1.012837479332056
$ws.Range("I17:N17").Value = $row17IN

# Row 18
$row18BF = New-Object "object[,]" 1,5
$row18BF[0,0] = 1.02
$row18BF[0,1] = 1.020223041089897
$row18BF[0,2] = 1.025639858714987
$row18BF[0,3] = 1.030480178827473
$row18BF[0,4] = 1.039655144386511
$ws.Range("B18:F18").Value = $row18BF
$row18IN = New-Object "object[,]" 1,6
$row18IN[0,0] = 1.028577568090549
$row18IN[0,1] = 1.026432099107601
$row18IN[0,2] = 1.028999817431747
$row18IN[0,3] = 1.033823207483648
$row18IN[0,4] = 1.042966557499708
$row18IN[0,5] = 1.012879218253267
$ws.Range("I18:N18").Value = $row18IN

# Row 19
$row19BF = New-Object "object[,]" 1,5
$row19BF[0,0] = 1.02
$row19BF[0,1] = 1.020293951012121
$row19BF[0,2] = 1.025692472671809
$row19BF[0,3] = 1.030548669829492
$row19BF[0,4] = 1.039743128484491
$ws.Range("B19:F19").Value = $row19BF
$row19IN = New-Object "object[,]" 1,6
$row19IN[0,0] = 1.028589503040743
$row19IN[0,1] = 1.026475083889072
$row19IN[0,2] = 1.029037776610791
$row19IN[0,3] = 1.033877068079503
$row19IN[0,4] = 1.043039994057688
$row19IN[0,5] = 1.012893445924257
$ws.Range("I19:N19").Value = $row19IN

# Row 20
$row20BF = New-Object "object[,]" 1,5
$row20BF[0,0] = 1.02
$row20BF[0,1] = 1.019976827055404
$row20BF[0,2] = 1.025457124949941
$row20BF[0,3] = 1.030242400929551
$row20BF[0,4] = 1.039349720487999
$ws.Range("B20:F20").Value = $row20BF
$row20IN = New-Object "object[,]" 1,6
$row20IN[0,0] = 1.02853596499083
$row20IN[0,1] = 1.026282800179099
$row20IN[0,2] = 1.028867914487027
$row20IN[0,3] = 1.033636179447734
$row20IN[0,4] = 1.042711600290134
$row20IN[0,5] = 1.012829799769254
$ws.Range("I20:N20").Value = $row20IN

# Row 21
$row21BF = New-Object "object[,]" 1,5
$row21BF[0,0] = 1.02
$row21BF[0,1] = 1.018945862668026
$row21BF[0,2] = 1.024691189636461
$row21BF[0,3] = 1.029247379201247
$row21BF[0,4] = 1.038072077369828
$ws.Range("B21:F21").Value = $row21BF
$row21IN = New-Object "object[,]" 1,6
$row21IN[0,0] = 1.028359063875936
$row21IN[0,1] = 1.025656868939998
$row21IN[0,2] = 1.02831392941318
$row21IN[0,3] = 1.032852830511851
$row21IN[0,4] = 1.041644521075595
$row21IN[0,5] = 1.012622588612879
$ws.Range("I21:N21").Value = $row21IN

# Row 22
$row22BF = New-Object "object[,]" 1,5
$row22BF[0,0] = 1.02
$row22BF[0,1] = 1.01829768585935
$row22BF[0,2] = 1.024209011685985
$row22BF[0,3] = 1.028622301677421
$row22BF[0,4] = 1.037269820868396
$ws.Range("B22:F22").Value = $row22BF
$row22IN = New-Object "object[,]" 1,6
$row22IN[0,0] = 1.028245662930866
$row22IN[0,1] = 1.025262713160231
$row22IN[0,2] = 1.027964283207065
$row22IN[0,3] = 1.032360160875515
$row22IN[0,4] = 1.04097403881805
$row22IN[0,5] = 1.012492084866612
$ws.Range("I22:N22").Value = $row22IN

# Row 23
$row23BF = New-Object "object[,]" 1,5
$row23BF[0,0] = 1.02
$row23BF[0,1] = 1.018641293436776
$row23BF[0,2] = 1.024464680144871
$row23BF[0,3] = 1.028953616355555
$row23BF[0,4] = 1.037695012666035
$ws.Range("B23:F23").Value = $row23BF
$row23IN = New-Object "object[,]" 1,6
$row23IN[0,0] = 1.028305984966366
$row23IN[0,1] = 1.025471719970618
$row23IN[0,2] = 1.028149763370374
$row23IN[0,3] = 1.032621348072501
$row23IN[0,4] = 1.041329432928299
$row23IN[0,5] = 1.012561288308284
$ws.Range("I23:N23").Value = $row23IN

# Row 24
$row24BF = New-Object "object[,]" 1,5
$row24BF[0,0] = 1.02
$row24BF[0,1] = 1.019994111791468
$row24BF[0,2] = 1.025469955618045
$row24BF[0,3] = 1.030259091554565
$row24BF[0,4] = 1.039371158080092
$ws.Range("B24:F24").Value = $row24BF
$row24IN = New-Object "object[,]" 1,6
$row24IN[0,0] = 1.028538893819618
$row24IN[0,1] = 1.026293283636929
$row24IN[0,2] = 1.028877179435306
$row24IN[0,3] = 1.033649309850814
$row24IN[0,4] = 1.042729497308779
$row24IN[0,5] = 1.012833269909897
$ws.Range("I24:N24").Value = $row24IN

# Row 25
$row25BF = New-Object "object[,]" 1,5
$row25BF[0,0] = 1.02
$row25BF[0,1] = 1.021564095246489
$row25BF[0,2] = 1.026633846126664
$row25BF[0,3] = 1.03177630912133
$row25BF[0,4] = 1.041320779236889
$ws.Range("B25:F25").Value = $row25BF
$row25IN = New-Object "object[,]" 1,6
$row25IN[0,0] = 1.028799672421043
$row25IN[0,1] = 1.027243989457239
$row25IN[0,2] = 1.029715463432962
$row25IN[0,3] = 1.034841537646848
$row25IN[0,4] = 1.044356059
$row25IN[0,5] = 1.013147913806756
$ws.Range("I25:N25").Value = $row25IN
